$wb = $excel.ActiveWorkbook
foreach ($name in @("Tir_213_40R21","Tir_270_70R22","Tir_145_70R13","Tir_235_50R24","Tir_430_50R38")) {
    $ws = $wb.Worksheets.Item($name)
    $rng = $ws.Range("A1")
    $fc = $rng.FormatConditions
    $cond = $fc.Add(1, 3, '"class"')
    $cond.Interior.Color = 13431551
}
Write-Output "done"
